$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the hourly crypto price/volume snapshot with the latest values.
# A couple of rows (Kaspa/Polygon) also swapped rank order. Numeric-looking
# price strings are given a leading apostrophe so Excel keeps storing them
# as text (matching the rest of the sheet) instead of auto-converting them
# to numbers.

$ws.Range("D2").Value = '59.688.85'
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = '2.616.47'
$ws.Range("E3").Value = '  +1.11%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''538.60'
$ws.Range("E5").Value = '  +2.57%  '
$ws.Range("D6").Value = '''142.48'
$ws.Range("E6").Value = '  +2.19%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("D9").Value = '''6.59'
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("E11").Value = '  +1.25%  '
$ws.Range("E12").Value = '  -1.41%  '
$ws.Range("D13").Value = '3.078.08'
$ws.Range("E13").Value = '  +1.08%  '
$ws.Range("D14").Value = '59.610.44'
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").Value = '''20.79'
$ws.Range("E15").Value = '  +1.26%  '
$ws.Range("D16").Value = '2.686.30'
$ws.Range("E16").Value = '  +4.56%  '
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").Value = '''341.01'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '''4.36'
$ws.Range("E19").Value = '  +1.05%  '
$ws.Range("D20").Value = '''10.13'
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = '''6.35'
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = '''67.26'
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").Value = '''0.409'
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").Value = '''0.166'
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("D26").Value = '''0.998'
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = '''7.23'
$ws.Range("E27").Value = '  +2.28%  '
$ws.Range("D28").Value = '0.0₃0746'
$ws.Range("E28").Value = '  +2.88%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '''1.67'
$ws.Range("E30").Value = '  +4.75%  '
$ws.Range("D31").Value = '''5.82'
$ws.Range("E31").Value = '  -2.27%  '
$ws.Range("D32").Value = '''18.80'
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("D33").Value = '''150.79'
$ws.Range("E33").Value = '  +1.06%  '
$ws.Range("D34").Value = '''3.99'
$ws.Range("E34").Value = '  +0.29%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = '''0.834'
$ws.Range("E36").Value = '  +3.25%  '
$ws.Range("E37").Value = '  -1.68%  '
$ws.Range("D38").Value = '''0.827'
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("D39").Value = '''3.54'
$ws.Range("E39").Value = '  +0.24%  '
$ws.Range("D40").Value = '''277.66'
$ws.Range("E40").Value = '  +2.20%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("E44").Value = '  -0.45%  '
$ws.Range("D45").Value = '''0.0527'
$ws.Range("E45").Value = '  +2.43%  '
$ws.Range("D46").Value = '1.956.64'
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("E47").Value = '  +0.30%  '
$ws.Range("D48").Value = '''18.47'
$ws.Range("E48").Value = '  +0.87%  '
$ws.Range("D49").Value = '''4.52'
$ws.Range("E49").Value = '  +1.36%  '
$ws.Range("D50").Value = '''112.00'
$ws.Range("E50").Value = '  -2.47%  '
$ws.Range("E51").Value = '  +0.68%  '
